$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 109. This shifts the existing rows
# 109..142 down to 110..143, preserving all of their data/formatting
# (matches the weekly-refresh pattern seen in the diff: every row from the
# previous run moves down by one, and a brand-new row 143 appears carrying
# what used to be row 142's data).
$ws.Rows(109).Insert()

# Populate the newly inserted row 109 with this week's new record.
$ws.Range("A109").Value = 8
$ws.Range("B109").Value = "Terminal La Palmera de La Serena"
$ws.Range("C109").Value = "Coquimbo"
$ws.Range("D109").Value = 44736
$ws.Range("D109").Style = $ws.Range("D108").Style
$ws.Range("E109").Value = 4
$ws.Range("F109").Value = 100112001
$ws.Range("G109").Value = "Berenjena"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 360
$ws.Range("K109").Value = 9000
$ws.Range("L109").Value = 10000
$ws.Range("M109").Value = 9500
$ws.Range("N109").Value = "`$/caja 50 unidades"
$ws.Range("O109").Value = "Región de Arica y Parinacota"
$ws.Range("P109").Value = 190
$ws.Range("Q109").Value = 50
$ws.Range("R109").Value = "Hortaliza"
